$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers (J: created_customer_id, K: creation_date) ---
$ws.Range("J1").Value = "created_customer_id"
$ws.Range("K1").Value = "creation_date"

# --- New data rows for created_customer_id (J2:J10) ---
$ws.Range("J2").Value = "cust8388157"
$ws.Range("J3").Value = "cust7914541"
$ws.Range("J4").Value = "cust2523276"
$ws.Range("J5").Value = "cust3015161"
$ws.Range("J6").Value = "cust5651821"
$ws.Range("J7").Value = "cust7775609"
$ws.Range("J8").Value = "cust5575266"
$ws.Range("J9").Value = "cust3817227"
$ws.Range("J10").Value = "cust7300786"

# --- New data rows for creation_date (K2:K10) - same literal text value ---
$ws.Range("K2").Value = "Jun 24, 2019"
$ws.Range("K3").Value = "Jun 24, 2019"
$ws.Range("K4").Value = "Jun 24, 2019"
$ws.Range("K5").Value = "Jun 24, 2019"
$ws.Range("K6").Value = "Jun 24, 2019"
$ws.Range("K7").Value = "Jun 24, 2019"
$ws.Range("K8").Value = "Jun 24, 2019"
$ws.Range("K9").Value = "Jun 24, 2019"
$ws.Range("K10").Value = "Jun 24, 2019"

# --- Formatting: existing header row A1:I1 gets a yellow highlight fill ---
$ws.Range("A1:I1").Interior.Color = 65535

# --- Formatting: new header cells J1:K1 get bold font + orange highlight fill ---
$ws.Range("J1:K1").Font.Bold = $true
$ws.Range("J1:K1").Interior.Color = 49407

# --- Column widths for the two new columns ---
$ws.Columns.Item(10).ColumnWidth = 18
$ws.Columns.Item(11).ColumnWidth = 12.25

# --- Update selection to K1 (matches the saved selection state in the file) ---
$null = $ws.Range("K1").Select()
